# Remove the leading "index" column (column A) from every table sheet in the
# workbook. Each sheet had an auto-generated integer index column (styled
# like the header row) that is no longer part of the expected output; the
# remaining columns shift one place to the left.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A:A").Delete()
}
